$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(265, 1).Value = 44789
$ws.Cells.Item(265, 2).Value = "KA03MP9766"
$ws.Cells.Item(265, 3).Value = "FIGO"
$ws.Cells.Item(265, 4).Value = "PMS"
$ws.Cells.Item(265, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(265, 6).Value = 32107
$ws.Cells.Item(265, 7).Value = "CREDIT"

$ws.Cells.Item(266, 1).Value = 44789
$ws.Cells.Item(266, 2).Value = "MH01AX9080"
$ws.Cells.Item(266, 3).Value = "POLO"
$ws.Cells.Item(266, 4).Value = "PMS"
$ws.Cells.Item(266, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(266, 6).Value = 14229

$ws.Cells.Item(267, 1).Value = 44789
$ws.Cells.Item(267, 2).Value = "DL7CP8758"
$ws.Cells.Item(267, 3).Value = "I20"
$ws.Cells.Item(267, 4).Value = "RR"
$ws.Cells.Item(267, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(267, 6).Value = 6000
$ws.Cells.Item(267, 7).Value = "G PAY"

$ws.Cells.Item(268, 1).Value = 44790
$ws.Cells.Item(268, 2).Value = "TN04AC3193"
$ws.Cells.Item(268, 3).Value = "I10"
$ws.Cells.Item(268, 4).Value = "WIRING PROBLEM"
$ws.Cells.Item(268, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(268, 6).Value = 1000
$ws.Cells.Item(268, 7).Value = "P PAY"

$ws.Cells.Item(269, 1).Value = 44790
$ws.Cells.Item(269, 2).Value = "KA03MN8120"
$ws.Cells.Item(269, 3).Value = "VENTO"
$ws.Cells.Item(269, 4).Value = "PMS                                      WW"
$ws.Cells.Item(269, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(270, 1).Value = 44790
$ws.Cells.Item(270, 2).Value = "KA03MV0364"
$ws.Cells.Item(270, 3).Value = "SCRPIO"
$ws.Cells.Item(270, 4).Value = "PMS                                      WW"
$ws.Cells.Item(270, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(271, 1).Value = 44790
$ws.Cells.Item(271, 2).Value = "KA51MB1424"
$ws.Cells.Item(271, 3).Value = "FIGO"
$ws.Cells.Item(271, 4).Value = "BRAKE PAD CHANGE"
$ws.Cells.Item(271, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(271, 6).Value = 1972
$ws.Cells.Item(271, 7).Value = "G PAY"

$ws.Cells.Item(272, 1).Value = 44790
$ws.Cells.Item(272, 2).Value = "KA03MN9595"
$ws.Cells.Item(272, 3).Value = "FORTUNER"
$ws.Cells.Item(272, 4).Value = "GENERAL CHECKUP         WW"
$ws.Cells.Item(272, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(273, 1).Value = 44790
$ws.Cells.Item(273, 2).Value = "KA03MJ1724"
$ws.Cells.Item(273, 3).Value = "COROLLA"
$ws.Cells.Item(273, 4).Value = "BRAKE DISC & PAD CHANGE"
$ws.Cells.Item(273, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(273, 6).Value = 8000
$ws.Cells.Item(273, 7).Value = "CASH"

$ws.Cells.Item(274, 1).Value = 44790
$ws.Cells.Item(274, 2).Value = "KL07BW6057"
$ws.Cells.Item(274, 3).Value = "INDICA "
$ws.Cells.Item(274, 4).Value = "HEAD LIGHT BULB"
$ws.Cells.Item(274, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(274, 6).Value = 1000
$ws.Cells.Item(274, 7).Value = "P PAY"

$ws.Cells.Item(275, 1).Value = 44790
$ws.Cells.Item(275, 2).Value = "KA03NE1061"
$ws.Cells.Item(275, 3).Value = "ECOSPORT"
$ws.Cells.Item(275, 4).Value = "HANDLE PROBLEM"
$ws.Cells.Item(275, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(275, 6).Value = 100
$ws.Cells.Item(275, 7).Value = "P PAY"

$ws.Cells.Item(276, 1).Value = 44791
$ws.Cells.Item(276, 2).Value = "KA01MS4265"
$ws.Cells.Item(276, 3).Value = "I10"
$ws.Cells.Item(276, 4).Value = "PMS"
$ws.Cells.Item(276, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(276, 6).Value = 3507
$ws.Cells.Item(276, 7).Value = "G PAY"

$ws.Cells.Item(277, 1).Value = 44791
$ws.Cells.Item(277, 2).Value = "WB20Z5652"
$ws.Cells.Item(277, 3).Value = "FIGO"
$ws.Cells.Item(277, 4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(277, 5).Value = "WORK DONE"

$ws.Cells.Item(278, 1).Value = 44791
$ws.Cells.Item(278, 2).Value = "MH14CC1351"
$ws.Cells.Item(278, 3).Value = "SWIFT"
$ws.Cells.Item(278, 4).Value = "STARTING PROBLEM"
$ws.Cells.Item(278, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(279, 1).Value = 44791
$ws.Cells.Item(279, 2).Value = "KA53MZ9550 "
$ws.Cells.Item(279, 3).Value = "ECOSPORT"
$ws.Cells.Item(279, 4).Value = "MIRROR CHANGE"
$ws.Cells.Item(279, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(280, 1).Value = 44791
$ws.Cells.Item(280, 2).Value = "CKE9802"
$ws.Cells.Item(280, 3).Value = " M-800"
$ws.Cells.Item(280, 4).Value = "PMS"
$ws.Cells.Item(280, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(280, 6).Value = 2000
$ws.Cells.Item(280, 7).Value = "G PAY"

$ws.Cells.Item(281, 1).Value = 44791
$ws.Cells.Item(281, 2).Value = "KA01MF7441"
$ws.Cells.Item(281, 3).Value = "MANZA"
$ws.Cells.Item(281, 4).Value = "PMS                                      WW"
$ws.Cells.Item(281, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(281, 6).Value = 6866
$ws.Cells.Item(281, 7).Value = "CREDIT"

$ws.Cells.Item(282, 1).Value = 44791
$ws.Cells.Item(282, 2).Value = "KA 05 MY7902"
$ws.Cells.Item(282, 3).Value = "NEXON"
$ws.Cells.Item(282, 4).Value = "BODY SHOP"
$ws.Cells.Item(282, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(282, 6).Value = 8250
$ws.Cells.Item(282, 7).Value = "G PAY"

$ws.Cells.Item(283, 1).Value = 44791
$ws.Cells.Item(283, 2).Value = "KA53MJ0986"
$ws.Cells.Item(283, 3).Value = "CELERIO"
$ws.Cells.Item(283, 4).Value = "MESH FITTING"
$ws.Cells.Item(283, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(283, 6).Value = 1300
$ws.Cells.Item(283, 7).Value = "G PAY"

$excel.ActiveWindow.ScrollRow = 262
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G283").Select()

